# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.977.24"
$ws.Range("E2").Value = "  +4.49%  "
$ws.Range("D3").Value = "3.364.16"
$ws.Range("E3").Value = "  +4.75%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'562.55"
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("D6").Value = "'153.69"
$ws.Range("E6").Value = "  +5.76%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'7.54"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("D11").Value = "'0.437"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "3.941.17"
$ws.Range("E12").Value = "  +4.78%  "
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'27.25"
$ws.Range("D15").Value = "'0.0000182"
$ws.Range("E15").Value = "  +3.63%  "
$ws.Range("D16").Value = "62.958.69"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "3.355.90"
$ws.Range("E17").Value = "  +4.49%  "
$ws.Range("D18").Value = "'6.47"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("D20").Value = "'8.45"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'389.97"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "'0.542"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'70.41"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +5.37%  "
$ws.Range("D26").Value = "'8.87"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "0.0₃0976"
$ws.Range("E27").Value = "  +7.21%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +6.30%  "
$ws.Range("D30").Value = "'2.00"
$ws.Range("E30").Value = "  +4.19%  "
$ws.Range("E31").Value = "  +3.92%  "
$ws.Range("D32").Value = "'23.07"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "  +6.83%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "  +8.70%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'161.04"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("E37").Value = "  +11.63%  "
$ws.Range("D38").Value = "'27.16"
$ws.Range("E38").Value = "  +5.03%  "
$ws.Range("D39").Value = "'0.0747"
$ws.Range("E39").Value = "  +4.89%  "
$ws.Range("D40").Value = "2.841.13"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "'0.0313"
$ws.Range("E41").Value = "  +9.27%  "
$ws.Range("D42").Value = "'4.33"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("D44").Value = "'40.79"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("D46").Value = "'22.31"
$ws.Range("E46").Value = "  +7.51%  "
$ws.Range("D47").Value = "3.406.16"
$ws.Range("E47").Value = "  +4.68%  "
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "'6.33"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").Value = "'0.812"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "'283.29"
$ws.Range("E51").Value = "  +5.23%  "
